# "Minor changes and Interventions runned in large scale mode"
#
# The workbook's `main` sheet drives a set of derived impact metrics (on
# sheets S, Y, Z, VA and within `main` itself) from a single assumption
# cell, C18: "Percentage of the smallholders to be covered". Switching the
# model to "large scale" bumps that coverage input way up; every dependent
# formula recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main")

$ws.Activate()
$ws.Range("C18").Value = 0.7904485917127797

# Reflect the cell the author was last working in as the active selection
# on the frozen-pane view.
[void]$ws.Range("C18").Select()
